$d = $word.ActiveDocument

# 1. Candidate name: DHARGAWE SATYATA DILESH -> ABHISHEK BAPNA
$d.Content.Find.Execute("DHARGAWE SATYATA DILESH", $true, $false, $false, $false, $false, $true, 1, $false, "ABHISHEK BAPNA", 2)

# 2. MBA Registration Number: MBA/0009/60 -> MBA/0003/60
$d.Content.Find.Execute("MBA Registration Number: MBA/0009/60", $true, $false, $false, $false, $false, $true, 1, $false, "MBA Registration Number: MBA/0003/60", 2)

# 3. Area of graduation checkbox: Commerce checked -> Science checked
$d.Content.Find.Execute("Engineering ☐    Science ☐    Commerce ☑    Arts ☐    Others ☐", $true, $false, $false, $false, $false, $true, 1, $false, "Engineering ☐    Science ☑    Commerce ☐    Arts ☐    Others ☐", 2)

# 4. Date: 11 April 2025 -> 15 April 2025
$d.Content.Find.Execute("Date: 11 April 2025", $true, $false, $false, $false, $false, $true, 1, $false, "Date: 15 April 2025", 2)

# 5. Remarks line
$d.Content.Find.Execute("Remarks:  |  |  |  |  |  |  |  |  |  |  |  |  |  |  |  | ", $true, $false, $false, $false, $false, $true, 1, $false, "Remarks: Enclosure II will be submitted at August 2nd                ", 2)
